# Scheduled-runner update: refresh the market-derived columns
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# for the rows whose underlying market data changed, across all
# Disciple-of-the-Hand sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21: Book and a Hard Place
$ws.Cells.Item(21, 8).Value = 1097
$ws.Cells.Item(21, 9).Value = 1097
$ws.Cells.Item(21, 11).Value = 1097
$ws.Cells.Item(21, 13).Value = -629

# Row 23: There's Something about Bury
$ws.Cells.Item(23, 8).Value = 1097
$ws.Cells.Item(23, 9).Value = 1097
$ws.Cells.Item(23, 11).Value = 1097
$ws.Cells.Item(23, 13).Value = -863

# Row 131: Mindful Study
$ws.Cells.Item(131, 8).Value = 2111.8572
$ws.Cells.Item(131, 9).Value = 779
$ws.Cells.Item(131, 11).Value = 2337
$ws.Cells.Item(131, 13).Value = 2703

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Cells.Item(2, 8).Value = 3224.75
$ws.Cells.Item(2, 9).Value = 3466.5
$ws.Cells.Item(2, 11).Value = 3466.5
$ws.Cells.Item(2, 13).Value = -3353.5

# Row 43: They've Got Legs
$ws.Cells.Item(43, 8).Value = 46342
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 14).ClearContents()

# Row 45: Hollow Hallmarks
$ws.Cells.Item(45, 8).Value = 2908.8333
$ws.Cells.Item(45, 9).Value = 2290.6
$ws.Cells.Item(45, 10).Value = 6000
$ws.Cells.Item(45, 11).Value = 2290.6
$ws.Cells.Item(45, 12).Value = 6000
$ws.Cells.Item(45, 13).Value = -1913.6
$ws.Cells.Item(45, 14).Value = -6754

# Row 61: Dealing with the Tough Stuff
$ws.Cells.Item(61, 8).Value = 3236.1875
$ws.Cells.Item(61, 10).Value = 2850
$ws.Cells.Item(61, 12).Value = 2850
$ws.Cells.Item(61, 14).Value = -3274

# Row 116: No Scope
$ws.Cells.Item(116, 8).Value = 3224.75
$ws.Cells.Item(116, 9).Value = 3466.5
$ws.Cells.Item(116, 11).Value = 3466.5
$ws.Cells.Item(116, 13).Value = -1172.5

# Row 136: Metal with Mettle
$ws.Cells.Item(136, 8).Value = 3236.1875
$ws.Cells.Item(136, 10).Value = 2850
$ws.Cells.Item(136, 12).Value = 8550
$ws.Cells.Item(136, 14).Value = -13650

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Cells.Item(3, 8).Value = 3224.75
$ws.Cells.Item(3, 9).Value = 3466.5
$ws.Cells.Item(3, 11).Value = 3466.5
$ws.Cells.Item(3, 13).Value = -3352.5

# Row 82: Spirituality Inspector
$ws.Cells.Item(82, 8).Value = 20030.9
$ws.Cells.Item(82, 9).Value = 9329.857
$ws.Cells.Item(82, 11).Value = 9329.857
$ws.Cells.Item(82, 13).Value = -8946.857

# Row 85: The Clamor for Hammers (L)
$ws.Cells.Item(85, 8).Value = 20030.9
$ws.Cells.Item(85, 9).Value = 9329.857
$ws.Cells.Item(85, 11).Value = 9329.857
$ws.Cells.Item(85, 13).Value = -8003.857

# Row 86: Through Thick and Thin
$ws.Cells.Item(86, 8).Value = 1132
$ws.Cells.Item(86, 9).Value = 1132
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 1132
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = -9
$ws.Cells.Item(86, 14).ClearContents()

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Cells.Item(89, 8).Value = 1132
$ws.Cells.Item(89, 9).Value = 1132
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 5660
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Value = -44
$ws.Cells.Item(89, 14).ClearContents()

# Row 97: File under Dull
$ws.Cells.Item(97, 8).Value = 28000
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 82: Aim to Please
$ws.Cells.Item(82, 8).Value = 22082
$ws.Cells.Item(82, 9).Value = 14164
$ws.Cells.Item(82, 10).Value = 30000
$ws.Cells.Item(82, 11).Value = 14164
$ws.Cells.Item(82, 12).Value = 30000
$ws.Cells.Item(82, 13).Value = -13803
$ws.Cells.Item(82, 14).Value = -30722

# Row 85: To Protect My City, I Must Wear a Mask (L)
$ws.Cells.Item(85, 8).Value = 22082
$ws.Cells.Item(85, 9).Value = 14164
$ws.Cells.Item(85, 10).Value = 30000
$ws.Cells.Item(85, 11).Value = 14164
$ws.Cells.Item(85, 12).Value = 30000
$ws.Cells.Item(85, 13).Value = -12916
$ws.Cells.Item(85, 14).Value = -32496

# Row 107: Built to Last
$ws.Cells.Item(107, 8).Value = 9000
$ws.Cells.Item(107, 9).Value = 1000
$ws.Cells.Item(107, 10).Value = 13000
$ws.Cells.Item(107, 11).Value = 1000
$ws.Cells.Item(107, 12).Value = 13000
$ws.Cells.Item(107, 13).Value = 920
$ws.Cells.Item(107, 14).Value = -16840

# Row 122: Timber of Tenkonto
$ws.Cells.Item(122, 8).Value = 861.5
$ws.Cells.Item(122, 10).Value = 874.75
$ws.Cells.Item(122, 12).Value = 2624.25
$ws.Cells.Item(122, 14).Value = -7524.25

$ws = $wb.Worksheets.Item("GSM")
# Row 120: A Beneficent Elegy
$ws.Cells.Item(120, 8).Value = 44500
$ws.Cells.Item(120, 10).Value = 44500
$ws.Cells.Item(120, 12).Value = 44500
$ws.Cells.Item(120, 14).Value = -54176

# Row 122: Awarding Academic Excellence
$ws.Cells.Item(122, 8).Value = 4167
$ws.Cells.Item(122, 9).Value = 4416
$ws.Cells.Item(122, 10).Value = 3669
$ws.Cells.Item(122, 11).Value = 13248
$ws.Cells.Item(122, 12).Value = 11007
$ws.Cells.Item(122, 13).Value = -10798
$ws.Cells.Item(122, 14).Value = -15907

# Row 126: Gold Rush Order
$ws.Cells.Item(126, 8).Value = 9796.6
$ws.Cells.Item(126, 9).Value = 9746.25
$ws.Cells.Item(126, 11).Value = 29238.75
$ws.Cells.Item(126, 13).Value = -26768.75

# Row 128: To Fight at Her Side
$ws.Cells.Item(128, 8).Value = 42000
$ws.Cells.Item(128, 10).Value = 42000
$ws.Cells.Item(128, 12).Value = 42000
$ws.Cells.Item(128, 14).Value = -51960

$ws = $wb.Worksheets.Item("LTW")
# Row 46: Supply Side Logic
$ws.Cells.Item(46, 8).Value = 600.3333
$ws.Cells.Item(46, 9).Value = 400.5
$ws.Cells.Item(46, 10).Value = 1000
$ws.Cells.Item(46, 11).Value = 400.5
$ws.Cells.Item(46, 12).Value = 1000
$ws.Cells.Item(46, 13).Value = -212.5
$ws.Cells.Item(46, 14).Value = -1376

# Row 48: Through a Glass Brightly
$ws.Cells.Item(48, 8).Value = 77500
$ws.Cells.Item(48, 9).Value = 77500
$ws.Cells.Item(48, 11).Value = 77500
$ws.Cells.Item(48, 13).Value = -76839

# Row 56: Hold On Tight
$ws.Cells.Item(56, 8).Value = 8582722
$ws.Cells.Item(56, 9).Value = 8582722
$ws.Cells.Item(56, 11).Value = 8582722
$ws.Cells.Item(56, 13).Value = -8582031

# Row 132: Tenets of Tanning
$ws.Cells.Item(132, 8).Value = 2822.111
$ws.Cells.Item(132, 9).Value = 2822.111
$ws.Cells.Item(132, 11).Value = 8466.332999999999
$ws.Cells.Item(132, 13).Value = -5936.332999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 2: The Unmentionables
$ws.Cells.Item(2, 8).Value = 6614.143
$ws.Cells.Item(2, 9).Value = 149.5
$ws.Cells.Item(2, 10).Value = 9200
$ws.Cells.Item(2, 11).Value = 149.5
$ws.Cells.Item(2, 12).Value = 9200
$ws.Cells.Item(2, 13).Value = -37.5
$ws.Cells.Item(2, 14).Value = -9424

# Row 4: Not Cool Enough
$ws.Cells.Item(4, 8).Value = 7898.4
$ws.Cells.Item(4, 10).Value = 11833.333
$ws.Cells.Item(4, 12).Value = 11833.333
$ws.Cells.Item(4, 14).Value = -12059.333

# Row 14: Hat in Hand
$ws.Cells.Item(14, 8).Value = 10000
$ws.Cells.Item(14, 9).Value = 10000
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 10000
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = -9832
$ws.Cells.Item(14, 14).ClearContents()

# Row 101: Who War It Better
$ws.Cells.Item(101, 8).Value = 32000
$ws.Cells.Item(101, 10).Value = 32000
$ws.Cells.Item(101, 12).Value = 32000
$ws.Cells.Item(101, 14).Value = -38490

# Row 107: Flax Wax
$ws.Cells.Item(107, 8).Value = 194.66667
$ws.Cells.Item(107, 9).Value = 184.5
$ws.Cells.Item(107, 11).Value = 553.5
$ws.Cells.Item(107, 13).Value = 1366.5

# Row 126: A Polished Purchase
$ws.Cells.Item(126, 8).Value = 4368.6665
$ws.Cells.Item(126, 9).Value = 5000
$ws.Cells.Item(126, 10).Value = 4053
$ws.Cells.Item(126, 11).Value = 15000
$ws.Cells.Item(126, 12).Value = 12159
$ws.Cells.Item(126, 13).Value = -12530
$ws.Cells.Item(126, 14).Value = -17099

